$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "코칩"
$ws.Range("C2").Value = "코스닥"
$ws.Range("D2").Value = 270
$ws.Range("E2").Value = "한국"
$ws.Range("F2").Value = 270
$ws.Range("G2").Value = "-"
$ws.Range("H2").Value = "-"
$ws.Range("I2").Value = "-"
$ws.Range("J2").Value = "-"
$ws.Range("K2").Value = "대표"
$ws.Range("L2").Value = "-"
$ws.Range("M2").Value = 18000
$ws.Range("N2").Value = 100
$ws.Range("Q2").Value = 1105000
$ws.Range("A2").Formula = '="2024-05-07"'
$ws.Range("O2").Formula = '="2024-04-24"'
$ws.Range("P2").Formula = '="2024-04-29"'
$ws.Range("A2:Q2").Copy()
$ws.Range("A2:Q2").PasteSpecial(-4163)

$ws.Range("B3").Value = "SK증권제12호스팩"
$ws.Range("C3").Value = "코스닥"
$ws.Range("D3").Value = 60
$ws.Range("E3").Value = "SK"
$ws.Range("F3").Value = 60
$ws.Range("G3").Value = "-"
$ws.Range("H3").Value = "-"
$ws.Range("I3").Value = "-"
$ws.Range("J3").Value = "-"
$ws.Range("K3").Value = "대표"
$ws.Range("L3").Value = "-"
$ws.Range("M3").Value = 2000
$ws.Range("N3").Value = 100
$ws.Range("Q3").Value = 2250000
$ws.Range("A3").Formula = '="2024-05-07"'
$ws.Range("O3").Formula = '="2024-04-23"'
$ws.Range("P3").Formula = '="2024-04-26"'
$ws.Range("A3:Q3").Copy()
$ws.Range("A3:Q3").PasteSpecial(-4163)

$ws.Range("B4").Value = "민테크"
$ws.Range("C4").Value = "코스닥"
$ws.Range("D4").Value = 315
$ws.Range("E4").Value = "KB"
$ws.Range("F4").Value = 315
$ws.Range("G4").Value = "-"
$ws.Range("H4").Value = "-"
$ws.Range("I4").Value = "-"
$ws.Range("J4").Value = "-"
$ws.Range("K4").Value = "대표"
$ws.Range("L4").Value = "-"
$ws.Range("M4").Value = 10500
$ws.Range("N4").Value = 100
$ws.Range("Q4").Value = 1800000
$ws.Range("A4").Formula = '="2024-05-03"'
$ws.Range("O4").Formula = '="2024-04-23"'
$ws.Range("P4").Formula = '="2024-04-26"'
$ws.Range("A4:Q4").Copy()
$ws.Range("A4:Q4").PasteSpecial(-4163)

$ws.Range("B5").Value = "디앤디파마텍"
$ws.Range("C5").Value = "코스닥"
$ws.Range("D5").Value = 363
$ws.Range("E5").Value = "한국"
$ws.Range("F5").Value = 363
$ws.Range("G5").Value = "-"
$ws.Range("H5").Value = "-"
$ws.Range("I5").Value = "-"
$ws.Range("J5").Value = "-"
$ws.Range("K5").Value = "대표"
$ws.Range("L5").Value = "-"
$ws.Range("M5").Value = 33000
$ws.Range("N5").Value = 100
$ws.Range("Q5").Value = 805400
$ws.Range("A5").Formula = '="2024-05-02"'
$ws.Range("O5").Formula = '="2024-04-22"'
$ws.Range("P5").Formula = '="2024-04-25"'
$ws.Range("A5:Q5").Copy()
$ws.Range("A5:Q5").PasteSpecial(-4163)

$ws.Range("B6").Value = "유안타제16호스팩"
$ws.Range("C6").Value = "코스닥"
$ws.Range("D6").Value = 103
$ws.Range("E6").Value = "유안타"
$ws.Range("F6").Value = 103
$ws.Range("G6").Value = "-"
$ws.Range("H6").Value = "-"
$ws.Range("I6").Value = "-"
$ws.Range("J6").Value = "-"
$ws.Range("K6").Value = "대표"
$ws.Range("L6").Value = "-"
$ws.Range("M6").Value = 2000
$ws.Range("N6").Value = 100
$ws.Range("Q6").Value = 3862500
$ws.Range("A6").Formula = '="2024-05-02"'
$ws.Range("O6").Formula = '="2024-04-22"'
$ws.Range("P6").Formula = '="2024-04-25"'
$ws.Range("A6:Q6").Copy()
$ws.Range("A6:Q6").PasteSpecial(-4163)

$ws.Range("B7").Value = "제일엠앤에스"
$ws.Range("C7").Value = "코스닥"
$ws.Range("D7").Value = 528
$ws.Range("E7").Value = "KB"
$ws.Range("F7").Value = 528
$ws.Range("G7").Value = "-"
$ws.Range("H7").Value = "-"
$ws.Range("I7").Value = "-"
$ws.Range("J7").Value = "-"
$ws.Range("K7").Value = "대표"
$ws.Range("L7").Value = "-"
$ws.Range("M7").Value = 22000
$ws.Range("N7").Value = 100
$ws.Range("Q7").Value = 1800000
$ws.Range("A7").Formula = '="2024-04-30"'
$ws.Range("O7").Formula = '="2024-04-18"'
$ws.Range("P7").Formula = '="2024-04-23"'
$ws.Range("A7:Q7").Copy()
$ws.Range("A7:Q7").PasteSpecial(-4163)

$ws.Range("B8").Value = "하나33호스팩"
$ws.Range("C8").Value = "코스닥"
$ws.Range("D8").Value = 70
$ws.Range("E8").Value = "하나"
$ws.Range("F8").Value = 70
$ws.Range("G8").Value = "-"
$ws.Range("H8").Value = "-"
$ws.Range("I8").Value = "-"
$ws.Range("J8").Value = "-"
$ws.Range("K8").Value = "대표"
$ws.Range("L8").Value = "-"
$ws.Range("M8").Value = 2000
$ws.Range("N8").Value = 100
$ws.Range("Q8").Value = 2625000
$ws.Range("A8").Formula = '="2024-04-24"'
$ws.Range("O8").Formula = '="2024-04-15"'
$ws.Range("P8").Formula = '="2024-04-18"'
$ws.Range("A8:Q8").Copy()
$ws.Range("A8:Q8").PasteSpecial(-4163)

$ws.Range("B9").Value = "신한제13호스팩"
$ws.Range("C9").Value = "코스닥"
$ws.Range("D9").Value = 60
$ws.Range("E9").Value = "신한"
$ws.Range("F9").Value = 60
$ws.Range("G9").Value = "-"
$ws.Range("H9").Value = "-"
$ws.Range("I9").Value = "-"
$ws.Range("J9").Value = "-"
$ws.Range("K9").Value = "대표"
$ws.Range("L9").Value = "-"
$ws.Range("M9").Value = 2000
$ws.Range("N9").Value = 100
$ws.Range("Q9").Value = 2250000
$ws.Range("A9").Formula = '="2024-04-22"'
$ws.Range("O9").Formula = '="2024-04-11"'
$ws.Range("P9").Formula = '="2024-04-15"'
$ws.Range("A9:Q9").Copy()
$ws.Range("A9:Q9").PasteSpecial(-4163)

$ws.Range("B10").Value = "신한제12호스팩"
$ws.Range("C10").Value = "코스닥"
$ws.Range("D10").Value = 100
$ws.Range("E10").Value = "신한"
$ws.Range("F10").Value = 100
$ws.Range("G10").Value = "-"
$ws.Range("H10").Value = "-"
$ws.Range("I10").Value = "-"
$ws.Range("J10").Value = "-"
$ws.Range("K10").Value = "대표"
$ws.Range("L10").Value = "-"
$ws.Range("M10").Value = 2000
$ws.Range("N10").Value = 100
$ws.Range("Q10").Value = 3750000
$ws.Range("A10").Formula = '="2024-04-15"'
$ws.Range("O10").Formula = '="2024-04-02"'
$ws.Range("P10").Formula = '="2024-04-05"'
$ws.Range("A10:Q10").Copy()
$ws.Range("A10:Q10").PasteSpecial(-4163)

$ws.Range("B11").Value = "아이엠비디엑스"
$ws.Range("C11").Value = "코스닥"
$ws.Range("D11").Value = 325
$ws.Range("E11").Value = "미래"
$ws.Range("F11").Value = 325
$ws.Range("G11").Value = "-"
$ws.Range("H11").Value = "-"
$ws.Range("I11").Value = "-"
$ws.Range("J11").Value = "-"
$ws.Range("K11").Value = "대표"
$ws.Range("L11").Value = "-"
$ws.Range("M11").Value = 13000
$ws.Range("N11").Value = 100
$ws.Range("Q11").Value = 1875000
$ws.Range("A11").Formula = '="2024-04-03"'
$ws.Range("O11").Formula = '="2024-03-25"'
$ws.Range("P11").Formula = '="2024-03-28"'
$ws.Range("A11:Q11").Copy()
$ws.Range("A11:Q11").PasteSpecial(-4163)

$ws.Range("B12").Value = "하나32호스팩"
$ws.Range("C12").Value = "코스닥"
$ws.Range("D12").Value = 60
$ws.Range("E12").Value = "하나"
$ws.Range("F12").Value = 60
$ws.Range("G12").Value = "-"
$ws.Range("H12").Value = "-"
$ws.Range("I12").Value = "-"
$ws.Range("J12").Value = "-"
$ws.Range("K12").Value = "대표"
$ws.Range("L12").Value = "-"
$ws.Range("M12").Value = 2000
$ws.Range("N12").Value = 100
$ws.Range("Q12").Value = 2250000
$ws.Range("A12").Formula = '="2024-03-27"'
$ws.Range("O12").Formula = '="2024-03-18"'
$ws.Range("P12").Formula = '="2024-03-21"'
$ws.Range("A12:Q12").Copy()
$ws.Range("A12:Q12").PasteSpecial(-4163)

$ws.Range("B13").Value = "엔젤로보틱스"
$ws.Range("C13").Value = "코스닥"
$ws.Range("D13").Value = 320
$ws.Range("E13").Value = "NH"
$ws.Range("F13").Value = 320
$ws.Range("G13").Value = "-"
$ws.Range("H13").Value = "-"
$ws.Range("I13").Value = "-"
$ws.Range("J13").Value = "-"
$ws.Range("K13").Value = "대표"
$ws.Range("L13").Value = "-"
$ws.Range("M13").Value = 20000
$ws.Range("N13").Value = 100
$ws.Range("Q13").Value = 880000
$ws.Range("A13").Formula = '="2024-03-26"'
$ws.Range("O13").Formula = '="2024-03-14"'
$ws.Range("P13").Formula = '="2024-03-19"'
$ws.Range("A13:Q13").Copy()
$ws.Range("A13:Q13").PasteSpecial(-4163)

$ws.Range("B14").Value = "삼현"
$ws.Range("C14").Value = "코스닥"
$ws.Range("D14").Value = 600
$ws.Range("E14").Value = "한국"
$ws.Range("F14").Value = 600
$ws.Range("G14").Value = "-"
$ws.Range("H14").Value = "-"
$ws.Range("I14").Value = "-"
$ws.Range("J14").Value = "-"
$ws.Range("K14").Value = "대표"
$ws.Range("L14").Value = "-"
$ws.Range("M14").Value = 30000
$ws.Range("N14").Value = 100
$ws.Range("Q14").Value = 1368000
$ws.Range("A14").Formula = '="2024-03-21"'
$ws.Range("O14").Formula = '="2024-03-12"'
$ws.Range("P14").Formula = '="2024-03-15"'
$ws.Range("A14:Q14").Copy()
$ws.Range("A14:Q14").PasteSpecial(-4163)

$excel.CutCopyMode = 0